$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Style = "Normal"
$ws.Range("A2").Value = "rua varginha"
$ws.Range("B2").Value = 38400322
$ws.Range("C2").Value = "CEP"

$ws.Range("B3").Select()
